$d = $word.ActiveDocument

# --- Insert "Best=[408263175]" as a new paragraph before the first paragraph ---
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphBefore()
$bestPara = $d.Paragraphs(1)
$bestPara.Range.Text = "Best=[408263175]"

# --- Insert "Worst=[351782064]" as a new paragraph right after it (still before the original first paragraph) ---
$origFirst = $d.Paragraphs(2)
$origFirst.Range.InsertParagraphBefore()
$worstPara = $d.Paragraphs(2)
$worstPara.Range.Text = "Worst=[351782064]"

# Re-seat the hidden "_GoBack" bookmark between "Worst=" and "[351782064]"
# (adding a bookmark with this reserved name replaces whichever bookmark
# previously held it, so the old one elsewhere in the doc disappears).
$worstStart = $worstPara.Range.Start
$bmPos = $worstStart + 6
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Clean up the (now third) paragraph: collapse the split runs/proofErr markers ---
$d.Content.Find.Execute("My ordering of move_expand_order is 4,0,2,6,8,1,3,5,7 as shown below:", $false, $false, $false, $false, $false, $true, 1, $false, "My ordering of move_expand_order is 4,0,2,6,8,1,3,5,7 as shown below:", 2) | Out-Null

# --- Merge the "42126" / " as shown below:" runs in the "Expanded node count" paragraph ---
$d.Content.Find.Execute("42126 as shown below:", $false, $false, $false, $false, $false, $true, 1, $false, "42126 as shown below:", 2) | Out-Null

# --- Merge the many split runs in the long explanation paragraph into one run ---
$oldExplanation = "order can get the count of expanded states less than 50,000 is explained as follows. I moved the positions, which is more likely to get a larger utility value, in front of the positions less likely to achieve this. Position 4 is the most likely position to win because it can be as one element for 4 possible cases (one row, one column and two diagonals). Then, position 0, 2, 6 and 8 can be as one element for 3 possible cases (one row, one column and one diagonals). At last, 1, 3, 5 and 7 can only make 2 possible cases (one row and one column). So, in this order, the algorithm can quickly get a larger alpha, which is more likely to "
$d.Content.Find.Execute($oldExplanation, $false, $false, $false, $false, $false, $true, 1, $false, $oldExplanation, 2) | Out-Null
